$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new practice-day column BD ---
# BD1: header date (2025-10-01 -> serial 45931), copy BC1's format (date style)
# then overwrite with the new date value.
$ws.Range("BC1").Copy($ws.Range("BD1"))
$ws.Range("BD1").Value = 45931

# For every player row that already has an entry in BC (the previous last
# day column), copy that same attendance marker into the new BD column.
# Row 12 has no BC cell (that player has no record for this period), so it
# is skipped, matching the source data.
$rows = (2..11) + (13..29)
foreach ($r in $rows) {
    $src = $ws.Range("BC$r")
    if ($src.Value -ne $null) {
        $src.Copy($ws.Range("BD$r"))
    }
}

# --- Update the view state to reflect the new rightmost data column ---
$ws.Range("BF23").Select()
$win = $excel.ActiveWindow
$win.ScrollColumn = 55
